$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 28
$wsMeans.Range("D9").Value = 44
$wsMeans.Range("E9").Value = 48
$wsMeans.Range("F9").Value = 54
$wsMeans.Range("G9").Value = 73

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.3
$wsMeans.Range("D10").Value = 0.3
$wsMeans.Range("E10").Value = 0.3
$wsMeans.Range("F10").Value = 0.3
$wsMeans.Range("G10").Value = 0.31

# --- Sheet "Standard Deviations" ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million)
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 13
$wsSD.Range("D9").Value = 8.2
$wsSD.Range("E9").Value = 5.4
$wsSD.Range("F9").Value = 13
$wsSD.Range("G9").Value = 51

# Row 10: Total Respiratory (hazard quotient)
$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.11
$wsSD.Range("D10").Value = 0
$wsSD.Range("E10").Value = 0
$wsSD.Range("F10").Value = 0.0000000000000000079
$wsSD.Range("G10").Value = 0.032
